$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.501.69"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "3.161.51"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'593.16"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'147.00"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D8").Value = "3.151.35"
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").Value = "'0.531"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("E11").Value = "  +3.90%  "
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").Value = "'37.39"
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").Value = "3.684.79"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").Value = "64.213.91"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").Value = "3.158.69"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "'469.38"
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("D21").Value = "'14.49"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = "  +8.34%  "
$ws.Range("D26").Value = "'81.55"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "'9.78"
$ws.Range("E28").Value = "  +9.56%  "
$ws.Range("E29").Value = "  +2.29%  "
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("D31").Value = "'7.40"
$ws.Range("E31").Value = "  +7.89%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "'27.65"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D35").Value = "0.0₃0845"
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("E37").Value = "  +3.30%  "
$ws.Range("D38").Value = "'2.32"
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("D39").Value = "'3.27"
$ws.Range("E39").Value = "  -3.47%  "
$ws.Range("D40").Value = "'471.11"
$ws.Range("E40").Value = "  +5.23%  "
$ws.Range("D41").Value = "'51.75"
$ws.Range("E42").Value = "  +6.09%  "
$ws.Range("D43").Value = "'0.297"
$ws.Range("E43").Value = "  +5.90%  "
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").Value = "2.939.57"
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("D46").Value = "'40.64"
$ws.Range("E46").Value = "  +11.67%  "
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").Value = "'129.70"
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").Value = "'2.27"
$ws.Range("E50").Value = "  +3.56%  "
$ws.Range("E51").Value = "  -0.15%  "
